$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.142.04'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.639.61'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.72'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.93'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '1.868.90'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.622.68'
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.541'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.82'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '27.161.30'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.24'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("E22").Value = '  +4.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.41'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.89'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.42'
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.67'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("D34").Value = '1.309.50'
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.545'
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  +5.77%  '
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("D44").Value = '1.779.27'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.78'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.46'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.66'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0962'
$ws.Range("E51").Value = '  -0.09%  '
